$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.703768610954285
$ws.Range("B1").Value = 5.387849807739258
$ws.Range("C1").Value = 6.046513080596924
$ws.Range("D1").Value = 2.574524641036987
$ws.Range("E1").Value = 1.739579677581787
